# Refactored code to make better use of Config constant, made Process more extensible
#
# Constants sheet:
#  - rename Status_Failed -> Status_Failure
#  - add a new Status_Pending / Pending row (right after Status_Failure)
#  - add EmailService_Gmail / Gmail and EmailService_Zoho / Zoho rows
#    (after the existing GoogleFormPrompt row, which shifts down by one row)

$wb = $excel.ActiveWorkbook

$wsSettings  = $wb.Worksheets.Item("Settings")
$wsConstants = $wb.Worksheets.Item("Constants")

# --- Constants sheet -------------------------------------------------

# 1. Rename the existing "Status_Failed" constant to "Status_Failure"
#    (value/meaning unchanged).
$wsConstants.Range("A14").Value = "Status_Failure"

# 2. Insert a new row below it for the new "Status_Pending" constant -
#    this pushes the blank row + GoogleFormPrompt block down by one.
$wsConstants.Rows.Item(15).Insert()
$wsConstants.Range("A15").Value = "Status_Pending"
$wsConstants.Range("B15").Value = "Pending"

# 3. Add the new email-service constants after the (now shifted)
#    GoogleFormPrompt row (row 17), leaving row 18 blank as a separator.
$wsConstants.Range("A19").Value = "EmailService_Gmail"
$wsConstants.Range("B19").Value = "Gmail"

$wsConstants.Range("A20").Value = "EmailService_Zoho"
$wsConstants.Range("B20").Value = "Zoho"

# Row height tweak that comes along with Excel re-flowing the wrapped
# header text when it resaved the file.
$wsConstants.Rows.Item(2).RowHeight = 28.8

# Restore the cursor to where the author left it.
$wsConstants.Activate()
$wsConstants.Range("A14").Select()

# --- Settings sheet ----------------------------------------------------

# Same wrapped-text row-height re-flow as above.
$wsSettings.Rows.Item(4).RowHeight = 28.8

"done"
